$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Shared-string / plain cell edits
# ---------------------------------------------------------------------------

# Version date in the intro cell
$ws.Range("A1").Value = "Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de ICTU Kwaliteitsaanpak Softwareontwikkeling versie wip, 22-09-2025."

# M32 sub-items: drop trailing comma/period/"en"
$ws.Range("B30").Value = "1. Een rapportage met tenminste de bevindingen, risico's voor opdrachtgevende organisatie en ICTU, en mitigerende maatregelen"
$ws.Range("B31").Value = "2. Een transitieplan dat de activiteiten beschrijft die nodig zijn om de software af te bouwen of te herbouwen en te onderhouden"
$ws.Range("B32").Value = "3. Als er significante technische schuld aanwezig is in de bestaande software: een plan voor het aflossen van deze schuld"

# M02 sub-items 1-5: trailing period -> comma (item 6 untouched)
$ws.Range("B34").Value = "1. Tijdens de voorfase: het project reviewt de deliverables periodiek"
$ws.Range("B35").Value = "2. Tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd de kwaliteit van de software"
$ws.Range("B36").Value = "3. Als operationeel beheer onderdeel is van de dienstverlening tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd het gedrag van de software in gebruik en beheer"
$ws.Range("B37").Value = "4. Tijdens de realisatiefase: het project evalueert periodiek en handmatig de kwaliteitseigenschappen van de software die niet geautomatiseerd kunnen worden gemeten"
$ws.Range("B38").Value = "5. Tijdens de realisatiefase: het project actualiseert en reviewt periodiek de documentatie"

# M07 sub-items: drop trailing comma/period
$ws.Range("B43").Value = "1. Bouw van de software"
$ws.Range("B44").Value = "2. Unit tests"
$ws.Range("B45").Value = "3. Regressietests"
$ws.Range("B46").Value = "4. Beveiligingstests"
$ws.Range("B47").Value = "5. Performancetests"
$ws.Range("B48").Value = "6. Toegankelijkheidstests"
$ws.Range("B49").Value = "7. Broncodekwaliteitscontroles"
$ws.Range("B50").Value = "8. Installatie van de software in test, acceptatie en/of productieomgevingen"
$ws.Range("B51").Value = "9. Oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgevende organisatie"

# M05 sub-items: drop trailing comma/period
$ws.Range("B77").Value = "1. Scrumteam bestaand uit product owner, ontwikkelaars (zoals programmeurs, testers en ontwerpers) en Scrummaster"
$ws.Range("B78").Value = "2. Proces met daily scrum, sprints, sprint planning, sprint review, sprint retrospective en sprint refinement"
$ws.Range("B79").Value = "3. Definition of Ready en Definition of Done"
$ws.Range("B80").Value = "4. Product backlog en sprint backlog"

# M34 sub-items: drop trailing comma/period and the paragraph-number parentheses
$ws.Range("B86").Value = "1. De documentatie beschrijft de ontwikkel- en testomgeving die is toegepast"
$ws.Range("B87").Value = "2. De functionele documentatie beschrijft gegevensmodellen, functionele indeling, koppelingen, berichtdefinities en workflows/processen"
$ws.Range("B88").Value = "3. Als operationeel beheer onderdeel was van de dienstverlening: de operationele bedieningsinstructies beschrijven minimaal back-up/recovery, procedures bij calamiteiten, regelmatig terugkerende beheeractiviteiten en opstart- en afsluitprocedures"
$ws.Range("B89").Value = "4. De product backlog bevat de bekende bugs en wensen"
$ws.Range("B90").Value = "5. De broncode kent een gezonde balans tussen isolatie, cohesie en koppeling"
$ws.Range("B91").Value = "6. De broncode heeft een beperkte mate van duplicatie"
$ws.Range("B92").Value = "7. De broncode heeft een beperkte mate van complexiteit"
$ws.Range("B93").Value = "8. De broncode bevat geen of een beperkt aantal niet-afgeronde werkzaamheden (`"todo's`")"
$ws.Range("B94").Value = "9. De tests raken een voldoende groot deel van de broncode"
$ws.Range("B95").Value = "10. De tests raken een voldoende groot deel van de functionaliteit (functionele dekking)"
$ws.Range("B96").Value = "11. De onderkende productrisico's zijn gedekt"
$ws.Range("B97").Value = "12. Er is een regressietest beschikbaar"
$ws.Range("B98").Value = "13. Er is traceerbaarheid van eisen naar testgevallen"
$ws.Range("B99").Value = "14. De testset is goed opgebouwd"

# ---------------------------------------------------------------------------
# 2) Comment (cell note) edits - in-session string replace, never printed,
#    to avoid console re-encoding of non-ASCII characters.
# ---------------------------------------------------------------------------

# B7 (M01): add a sentence about when an IAMA is used, right after the
# existing IAMA explanation paragraph.
$c = $ws.Range("B7").Comment
$t = $c.Text()
$old = "In het IAMA worden verbanden gelegd met relevante regels, instrumenten en toetskaders op het gebied van algoritmen.`n`nZie https://www.rijksoverheid.nl"
$new = "In het IAMA worden verbanden gelegd met relevante regels, instrumenten en toetskaders op het gebied van algoritmen.`n`nEen IAMA wordt ingezet in alle gevallen waarin een overheidsorgaan overweegt een algoritme te (laten) ontwikkelen, in te kopen, aan te passen en/of in te gaan zetten.`n`nZie https://www.rijksoverheid.nl"
$c.Text($t.Replace($old, $new))

# B33 (M02): trailing period -> comma for items 1-5 (item 6 unaffected)
$c = $ws.Range("B33").Comment
$t = $c.Text()
$old = "1. Tijdens de voorfase: het project reviewt de deliverables periodiek.`n2. Tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd de kwaliteit van de software.`n3. Als operationeel beheer onderdeel is van de dienstverlening tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd het gedrag van de software in gebruik en beheer.`n4. Tijdens de realisatiefase: het project evalueert periodiek en handmatig de kwaliteitseigenschappen van de software die niet geautomatiseerd kunnen worden gemeten.`n5. Tijdens de realisatiefase: het project actualiseert en reviewt periodiek de documentatie.`n6."
$new = "1. Tijdens de voorfase: het project reviewt de deliverables periodiek,`n2. Tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd de kwaliteit van de software,`n3. Als operationeel beheer onderdeel is van de dienstverlening tijdens de realisatiefase: het project bewaakt op dagelijkse basis en geautomatiseerd het gedrag van de software in gebruik en beheer,`n4. Tijdens de realisatiefase: het project evalueert periodiek en handmatig de kwaliteitseigenschappen van de software die niet geautomatiseerd kunnen worden gemeten,`n5. Tijdens de realisatiefase: het project actualiseert en reviewt periodiek de documentatie,`n6."
$c.Text($t.Replace($old, $new))

# B75 (M23): add "inclusief projectleider" + teamleden -> projectleden (x2)
$c = $ws.Range("B75").Comment
$t = $c.Text()
$old = "Projectleden die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak.`n`nRationale`n`nHet inzetten van teamleden die bekend zijn met de Kwaliteitsaanpak zorgt voor een soepeler start van een nieuw project omdat zij bekend zijn met de inhoud van de Kwaliteitsaanpak, zoals kwaliteitsnormen en tools, en omdat zij al doende nieuwe teamleden bekend kunnen maken met de Kwaliteitsaanpak."
$new = "Projectleden, inclusief projectleider, die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak.`n`nRationale`n`nHet inzetten van projectleden die bekend zijn met de Kwaliteitsaanpak zorgt voor een soepeler start van een nieuw project omdat zij bekend zijn met de inhoud van de Kwaliteitsaanpak, zoals kwaliteitsnormen en tools, en omdat zij al doende nieuwe projectleden bekend kunnen maken met de Kwaliteitsaanpak."
$c.Text($t.Replace($old, $new))

# B85 (M34): drop "(code dekking)" from sub-item 9
$c = $ws.Range("B85").Comment
$t = $c.Text()
$old = "9. De tests raken een voldoende groot deel van de broncode (code dekking) (7.1),"
$new = "9. De tests raken een voldoende groot deel van de broncode (7.1),"
$c.Text($t.Replace($old, $new))
